# Apply cell value changes to column A as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    25 = 0
    27 = 1
    44 = 1
    58 = 1
    65 = 1
    118 = 0
    121 = 0
    139 = 0
    155 = 0
    158 = 1
    169 = 1
    196 = 1
    212 = 1
    255 = 1
    259 = 0
    304 = 0
    333 = 1
    345 = 0
    351 = 0
    361 = 1
    392 = 0
    403 = 0
    431 = 1
    434 = 0
    495 = 0
    558 = 0
    581 = 1
    598 = 1
    607 = 0
    613 = 1
    614 = 1
    622 = 1
    637 = 0
    649 = 1
    650 = 0
    656 = 0
    662 = 1
    664 = 1
    681 = 1
    733 = 0
    735 = 0
    806 = 1
    859 = 0
    913 = 1
    916 = 0
    917 = 1
    924 = 1
    933 = 0
    936 = 1
    938 = 1
    961 = 0
    976 = 0
    978 = 1
    986 = 1
    987 = 0
    988 = 0
    996 = 0
    998 = 0
    1000 = 0
    1021 = 0
    1031 = 1
    1043 = 0
    1045 = 1
    1056 = 1
    1060 = 1
    1072 = 1
    1099 = 0
    1160 = 1
    1198 = 0
    1217 = 1
    1230 = 0
    1233 = 0
    1254 = 1
    1435 = 1
    1437 = 1
    1439 = 1
    1441 = 1
    1442 = 1
    1443 = 1
    1445 = 1
    1446 = 1
    1447 = 1
    1448 = 1
    1449 = 1
    1451 = 1
    1454 = 1
    1455 = 1
    1457 = 1
    1458 = 1
    1459 = 1
    1460 = 1
    1461 = 1
    1463 = 1
    1464 = 1
    1466 = 1
    1467 = 1
    1468 = 1
    1469 = 1
    1470 = 1
    1471 = 1
    1472 = 1
    1473 = 1
    1474 = 1
    1476 = 0
    1480 = 1
    1482 = 1
    1484 = 0
    1501 = 1
    1504 = 1
    1505 = 1
    1506 = 0
    1508 = 0
    1517 = 0
    1518 = 1
    1523 = 1
    1528 = 0
    1531 = 0
    1532 = 0
    1535 = 0
    1539 = 0
    1542 = 0
    1546 = 0
    1548 = 1
    1549 = 1
    1550 = 1
    1552 = 1
    1556 = 0
    1559 = 1
    1560 = 0
    1565 = 0
    1566 = 0
    1568 = 1
    1570 = 1
    1571 = 1
    1572 = 0
    1573 = 1
    1574 = 1
    1577 = 0
    1580 = 1
    1581 = 0
    1582 = 0
    1585 = 0
    1587 = 0
    1589 = 1
    1590 = 0
    1591 = 0
    1594 = 0
    1595 = 0
    1597 = 1
    1598 = 1
    1599 = 0
    1600 = 0
    1601 = 0
    1607 = 0
    1608 = 1
    1614 = 1
    1615 = 1
    1616 = 0
    1617 = 0
    1618 = 0
    1620 = 0
    1621 = 0
    1623 = 0
    1625 = 0
    1626 = 0
    1628 = 0
    1629 = 0
    1630 = 0
    1631 = 0
    1632 = 0
    1633 = 0
    1634 = 0
    1635 = 0
    1636 = 0
    1637 = 0
    1638 = 0
    1639 = 0
    1640 = 0
    1641 = 0
    1642 = 0
    1643 = 0
    1644 = 0
    1645 = 0
    1646 = 0
    1647 = 0
    1648 = 0
    1649 = 0
    1651 = 0
    1657 = 1
    1660 = 0
    1675 = 0
    1677 = 0
    1702 = 1
    1717 = 1
    1736 = 0
    1745 = 1
    1765 = 1
    1784 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 1).Value = $changes[$row]
}

